# Refresh crypto price/symbol listing (coinranking.com scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(CellAddress, NewValue). Values are written through a
# text ('@') number format so numeric-looking strings (e.g. "246.10")
# are stored verbatim - preserving trailing zeros / exact precision -
# instead of being parsed into a Double. The style is reset back to
# "Normal" immediately afterwards so cell formatting is left untouched.
$changes = @(
    @('D2', '246.10'),
    @('D3', '21.99'),
    @('D4', '5.374'),
    @('D5', '0.05799'),
    @('D6', '3.378'),
    @('D8', '0.8092'),
    @('D9', '1.014'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D10', '0.1426'),
    @('E10', '9WazirXWRX'),
    @('B11', 'MandalaExchangeToken'),
    @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D11', '0.07501'),
    @('E11', '10MandalaExchangeTokenMDX'),
    @('B12', 'LiechtensteinCryptoassetsExchange'),
    @('C12', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D12', '0.03205'),
    @('E12', '11LiechtensteinCryptoassetsExchangeLCX'),
    @('B13', 'BitrueCoin'),
    @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D13', '0.03031'),
    @('E13', '12BitrueCoinBTR'),
    @('B14', 'MCDex'),
    @('C14', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'),
    @('D14', '4.182'),
    @('E14', '13MCDexMCB'),
    @('B15', 'BitMartToken'),
    @('C15', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D15', '0.09393'),
    @('E15', '14BitMartTokenBMX'),
    @('B16', 'BitForexToken'),
    @('C16', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D16', '0.001596'),
    @('E16', '15BitForexTokenBF'),
    @('B17', 'CoinExToken'),
    @('C17', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
    @('D17', '0.04807'),
    @('E17', '16CoinExTokenCET'),
    @('B18', 'One'),
    @('C18', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D18', '0.0005900'),
    @('E18', '17OneONE'),
    @('D19', '0.005651'),
    @('D20', '0.004095'),
    @('E20', '19HotbitTokenHTB'),
    @('D21', '0.0009972'),
    @('D23', '3.704'),
    @('D24', '2.243'),
    @('D26', '0.1296'),
    @('D27', '0.0003998'),
    @('D40', '0.03879'),
    @('B41', 'BKEXToken'),
    @('C41', 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'),
    @('D41', '0.1073'),
    @('E41', '40BKEXTokenBKK'),
    @('B42', 'CEJI'),
    @('C42', 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'),
    @('D42', '0.002640'),
    @('E42', '41CEJICEJI'),
    @('B43', 'KickToken'),
    @('C43', 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'),
    @('D43', '0.003043'),
    @('E43', '42KickTokenKICKWorstin24h'),
    @('D44', '0.006687'),
    @('D45', '0.00005592'),
    @('D47', '0.3900'),
    @('D48', '0.1469'),
    @('D49', '0.00002100')
)

foreach ($change in $changes) {
    $addr = $change[0]
    $val = $change[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
